$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update accuracy values for Random forest classifier rows (B5:B8)
$ws.Range("B5").Value = 0.68018433179723503
$ws.Range("B6").Value = 0.69216589861751154
$ws.Range("B7").Value = 0.69216589861751154
$ws.Range("B8").Value = 0.68663594470046085

# Update the active selection to match the saved view state
$ws.Range("E6").Select()
